$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: 星期一 8:10~9:00 101 密碼學 陳老師 資電101 用程式寫密碼 30 30 40 / 3 / 60 / 0
$ws.Range("A12").Value = "星期一"
$ws.Range("B12").Value = "8:10~9:00"
$ws.Range("C12").Value = "101"
$ws.Range("D12").Value = "密碼學"
$ws.Range("E12").Value = "陳老師"
$ws.Range("F12").Value = "資電101"
$ws.Range("G12").Value = "用程式寫密碼"
$ws.Range("H12").Value = "30 30 40"
$ws.Range("I12").Value = 3
$ws.Range("J12").Value = 60
$ws.Range("K12").Value = 0

# Row 13: 星期一 9:10~10:00 102 離散數學 李老師 資電102 就是離散 30 30 40 / 3 / 60 / 5
$ws.Range("A13").Value = "星期一"
$ws.Range("B13").Value = "9:10~10:00"
$ws.Range("C13").Value = "102"
$ws.Range("D13").Value = "離散數學"
$ws.Range("E13").Value = "李老師"
$ws.Range("F13").Value = "資電102"
$ws.Range("G13").Value = "就是離散"
$ws.Range("H13").Value = "30 30 40"
$ws.Range("I13").Value = 3
$ws.Range("J13").Value = 60
$ws.Range("K13").Value = 5

# Row 14: 星期三 8:10~9:00 103 微積分 王老師 資電103 危機分 30 30 40 / 3 / 60 / 6
$ws.Range("A14").Value = "星期三"
$ws.Range("B14").Value = "8:10~9:00"
$ws.Range("C14").Value = "103"
$ws.Range("D14").Value = "微積分"
$ws.Range("E14").Value = "王老師"
$ws.Range("F14").Value = "資電103"
$ws.Range("G14").Value = "危機分"
$ws.Range("H14").Value = "30 30 40"
$ws.Range("I14").Value = 3
$ws.Range("J14").Value = 60
$ws.Range("K14").Value = 6

# Update the active selection / view to match the saved worksheet state
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("I16").Select()
